$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at B (old B shifts to D, old C shifts to E).
# This mirrors the diff: B1 gets a brand-new header, old B1 header moves to D1,
# old C1 header (and its data column) moves to E1, and two fresh columns
# (new B, new C) are created with blank ("UN") data cells.
$ws.Range("B:C").Insert()

# Restore/assign the explicit custom width (8.0 "chars" in the stored XML,
# which this engine reports back as ColumnWidth 7.1666667) on every column
# that should carry a customWidth flag: the two freshly inserted columns and
# the one that got shifted over from the original layout.
$ws.Columns("C:C").ColumnWidth = 7.1666667
$ws.Columns("D:D").ColumnWidth = 7.1666667
$ws.Columns("E:E").ColumnWidth = 7.1666667

# New headers for the inserted columns (set in sheet order so the new shared
# strings are appended Jun_15 then Jun_17, matching the source order).
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Fill data rows (2-27) in the two new columns with the same blank marker
# ("UN") that's already used throughout the sheet for "no data" cells.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}
